$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add two new bulleted ("ListParagraph", numId 14) entries right after the
#    "Password reset UI added" bullet and before the "Free Thought and TTD"
#    Heading1.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Password reset UI added*") {
        $rng = $p.Range
        $rng.Collapse(0)                      # wdCollapseEnd
        $rng.InsertParagraphAfter()
        $p1 = $p.Next()
        $p1.Range.Text = "Fixed multiple issue on landscape viewport "

        $rng2 = $p1.Range
        $rng2.Collapse(0)                     # wdCollapseEnd
        $rng2.InsertParagraphAfter()
        $p2 = $p1.Next()
        $p2.Range.Text = "Created temporary UI for cart "
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Mark the "Tested function" (Heading2) and "Success in receiving data"
#    (Heading4) runs as carrying a rendered page break, i.e. add a
#    <w:lastRenderedPageBreak/> as the very first child of their run, before
#    the <w:t>.
# ---------------------------------------------------------------------------
function Add-LastRenderedPageBreak([string]$searchText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$searchText*") {
            $rng = $p.Range
            $full = $rng.WordOpenXML
            if ($full -match '<w:p [^>]*>.*?</w:p>') {
                $paraXml = $matches[0]
                $needle = '<w:r><w:t>' + $searchText + '</w:t></w:r>'
                $replacement = '<w:r><w:lastRenderedPageBreak/><w:t>' + $searchText + '</w:t></w:r>'
                $newParaXml = $paraXml.Replace($needle, $replacement)
                $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
                       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                       '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
                       '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
                $rng.InsertXML($pkg)
            }
            break
        }
    }
}

Add-LastRenderedPageBreak "Tested function"
Add-LastRenderedPageBreak "Success in receiving data"
